# Executive Summary.docx edits
$d = $word.ActiveDocument

# 1) Remove the stray "_GoBack" bookmark that originally sat between
#    "April 14" and ", 2017" in the date line.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) "(time efficiency:" -> "(average case time efficiency:" (both
#    occurrences in the algorithm-selection paragraph).
$d.Content.Find.Execute(
    "(time efficiency:", $true, $false, $false, $false, $false,
    $true, 1, $false, "(average case time efficiency:", 2) | Out-Null

# 3) "we chose small, medium, and large dataset sizes." ->
#    "small, medium, and large dataset sizes were chosen."
$d.Content.Find.Execute(
    "we chose small, medium, and large dataset sizes.", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "small, medium, and large dataset sizes were chosen.", 2) | Out-Null

# 4) "I felt it showed more of a " -> "it showed a more profound "
$d.Content.Find.Execute(
    "I felt it showed more of a ", $true, $false, $false, $false, $false,
    $true, 1, $false, "it showed a more profound ", 2) | Out-Null

# 5) "...n2) efficiency." -> "...n2) time efficiency. "
$d.Content.Find.Execute(
    ") efficiency.", $true, $false, $false, $false, $false,
    $true, 1, $false, ") time efficiency. ", 2) | Out-Null

# 6) "used as divide and conquer, is O(n)" ->
#    "used as the divide and conquer algorithm, is O(n)"
$d.Content.Find.Execute(
    "used as divide and conquer, is O(n)", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "used as the divide and conquer algorithm, is O(n)", 2) | Out-Null

# 7) Re-add the "_GoBack" bookmark at its new location: right after
#    "...divide and conquer algorithm" and before ", is O(n)".
$r = $d.Content
$r.Find.Execute(
    "used as the divide and conquer algorithm", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
